$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Beast', ['Token Creature " + [char]8212 + " Beast', '3/3'])"
$ws.Range("A3").Value = "('Goblin', ['Token Creature " + [char]8212 + " Goblin', '1/1'])"
$ws.Range("A4").Value = "('Golem', ['Token Artifact Creature " + [char]8212 + " Golem', '3/3'])"
$ws.Range("A5").Value = "('Myr', ['Token Artifact Creature " + [char]8212 + " Myr', '1/1'])"
$ws.Range("A6").Value = "('Poison Counter', ['Card', '(A player with ten or more poison counters loses the game.)'])"

$ws.Range("A7:A16").EntireRow.Delete()
